# full_report.xlsx — "change report struct"
#   * "Итог" loses the "Результат" / "Результат (код)" columns, keeping only
#     Результат (описание) / Количество URL / Процент, and its title bar
#     gets a bottom border instead of none.
#   * "Детальный отчет" becomes the active tab, gains an AutoFilter on its
#     last column (and the hidden _FilterDatabase name that goes with it).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Итог"
$ws2 = $wb.Worksheets.Item(2)   # "Детальный отчет"

# ---------------------------------------------------------------------
# 1. "Итог": drop the "Результат" (col A) and "Результат (код)" (col C)
#    columns, leaving Результат (описание) / Количество URL / Процент.
# ---------------------------------------------------------------------
[void]$ws1.Columns.Item(3).Delete()   # old column C: "Результат (код)"
[void]$ws1.Columns.Item(1).Delete()   # old column A: "Результат"

# re-merge the (now 3-wide) title row and give it a bottom border
[void]$ws1.Range("A1:C1").Merge()
$ws1.Range("A1:C1").Borders.Item(9).LineStyle = 1

# sheet is no longer the active tab; selection moves to C15
[void]$ws1.Range("C15").Select()

# ---------------------------------------------------------------------
# 2. "Детальный отчет": turn on an AutoFilter for the last column and
#    make this sheet the active tab with a fresh selection.
# ---------------------------------------------------------------------
[void]$ws2.Range("G1:G2").AutoFilter()
$filterName = $ws2.Names.Add("_xlnm._FilterDatabase", "='Детальный отчет'!`$G`$1:`$G`$2")
$filterName.Visible = $false

[void]$ws2.Activate()
[void]$ws2.Range("E7").Select()

Write-Output "ok"
